$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: move skill_bonus_per_level (L) value into base_damage_mod_bonus_per_level (E), clear L ---
$ws.Range("E11").Value = 0.001
$ws.Range("L11").ClearContents()

# --- Row 12: same pattern, keep F12 untouched ---
$ws.Range("E12").Value = 0.001
$ws.Range("L12").ClearContents()

# --- Row 14: fight_time_out_mod_bonus_per_level (H) value fix ---
$ws.Range("H14").Value = 0.00005

# --- Row 16: fight_time_out_mod_bonus_per_level (H) value fix ---
$ws.Range("H16").Value = 0.00004

# --- Row 17: is_locked (N) flag fix ---
$ws.Range("N17").Value = 0

# --- New row 19: Fighters Resilience ---
$ws.Range("A19").Value = "Fighters Resilience"
$ws.Range("B19").Value = 5
$ws.Range("C19").Value = "The attack timer is 10 seconds by default, how ever with this skill, you can shave off - at max level roughly 20% of those 10 seconds, letting you click more. You level this skill via fighting monsters. Click train on the character sheet and assign some xp, the xp you assign is what you will sacrifice from killing monsters. The higher the percentage, the less xp you get towards leveling when fighting monsters."
$ws.Range("D19").Value = 999
$ws.Range("H19").Value = 0.0002
$ws.Range("K19").Value = 1
$ws.Range("N19").Value = 0

# --- New row 20: Quick Feet ---
$ws.Range("A20").Value = "Quick Feet"
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = "As you move around the map, using directional buttons only (this will not apply to teleporting, setting sail or teleporting to your own kingdom), you will find, as you get higher in level - that the time is reducing. Train this skill by fighting monsters, to do so - click train on the character sheet for this skill and select how much xp you want to sacrifice to train this skill. Special locations with adventures, can have bonuses applied to skill training."
$ws.Range("D20").Value = 999
$ws.Range("I20").Value = 0.0002
$ws.Range("K20").Value = 1
$ws.Range("N20").Value = 0

# --- Column width tweaks (account for the COM ColumnWidth <-> XML width padding offset) ---
$ws.Columns.Item(1).ColumnWidth = 22.17
$ws.Columns.Item(3).ColumnWidth = 543.17

# --- Update selection / view position to match author's final cursor position ---
[void]$ws.Range("H14").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 5
